$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CharChanges")

$ws.Range("A11").Value = "Sabine"
$ws.Range("B11").Value = "Start SLP 4 -> 15"

$ws.Range("A12").Value = "Sabine"
$ws.Range("B12").Value = "SLP/Lvl 10 -> 18"

$ws.Range("A13").Value = "Sabine"
$ws.Range("B13").Value = "Replace some start items"

$ws.Range("A14").Value = "Gryban"
$ws.Range("B14").Value = "SLP/Lvl 3 -> 5"

$ws.Range("A15").Value = "Gryban"
$ws.Range("B15").Value = "Start SLP 19 -> 20"

$ws.Range("A16").Select()

$wb.Save()
